$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$c = $ws.Cells.Item(2, 4)
$c.NumberFormat = "@"
$c.Value = '43.210.29'
$c.Style = "Normal"
$c = $ws.Cells.Item(2, 5)
$c.NumberFormat = "@"
$c.Value = '  +1.57%  '
$c.Style = "Normal"

# Row 3
$c = $ws.Cells.Item(3, 4)
$c.NumberFormat = "@"
$c.Value = '2.448.75'
$c.Style = "Normal"
$c = $ws.Cells.Item(3, 5)
$c.NumberFormat = "@"
$c.Value = '  -2.61%  '
$c.Style = "Normal"

# Row 4
$c = $ws.Cells.Item(4, 4)
$c.NumberFormat = "@"
$c.Value = '0.993'
$c.Style = "Normal"
$c = $ws.Cells.Item(4, 5)
$c.NumberFormat = "@"
$c.Value = '  -0.59%  '
$c.Style = "Normal"

# Row 5
$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = '315.88'
$c.Style = "Normal"
$c = $ws.Cells.Item(5, 5)
$c.NumberFormat = "@"
$c.Value = '  -0.23%  '
$c.Style = "Normal"

# Row 6
$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = "@"
$c.Value = '98.21'
$c.Style = "Normal"
$c = $ws.Cells.Item(6, 5)
$c.NumberFormat = "@"
$c.Value = '  +4.38%  '
$c.Style = "Normal"

# Row 7
$c = $ws.Cells.Item(7, 4)
$c.NumberFormat = "@"
$c.Value = '0.579'
$c.Style = "Normal"
$c = $ws.Cells.Item(7, 5)
$c.NumberFormat = "@"
$c.Value = '  +0.15%  '
$c.Style = "Normal"

# Row 8
$c = $ws.Cells.Item(8, 4)
$c.NumberFormat = "@"
$c.Value = '0.998'
$c.Style = "Normal"
$c = $ws.Cells.Item(8, 5)
$c.NumberFormat = "@"
$c.Value = '  -0.31%  '
$c.Style = "Normal"

# Row 9
$c = $ws.Cells.Item(9, 4)
$c.NumberFormat = "@"
$c.Value = '0.541'
$c.Style = "Normal"
$c = $ws.Cells.Item(9, 5)
$c.NumberFormat = "@"
$c.Value = '  +2.34%  '
$c.Style = "Normal"

# Row 10
$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = "@"
$c.Value = '36.03'
$c.Style = "Normal"
$c = $ws.Cells.Item(10, 5)
$c.NumberFormat = "@"
$c.Value = '  +1.47%  '
$c.Style = "Normal"

# Row 11
$c = $ws.Cells.Item(11, 5)
$c.NumberFormat = "@"
$c.Value = '  +0.47%  '
$c.Style = "Normal"

# Row 12
$c = $ws.Cells.Item(12, 5)
$c.NumberFormat = "@"
$c.Value = '  +0.63%  '
$c.Style = "Normal"

# Row 13
$c = $ws.Cells.Item(13, 4)
$c.NumberFormat = "@"
$c.Value = '2.993.34'
$c.Style = "Normal"
$c = $ws.Cells.Item(13, 5)
$c.NumberFormat = "@"
$c.Value = '  +3.27%  '
$c.Style = "Normal"

# Row 14
$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = "@"
$c.Value = '0.108'
$c.Style = "Normal"
$c = $ws.Cells.Item(14, 5)
$c.NumberFormat = "@"
$c.Value = '  -0.17%  '
$c.Style = "Normal"

# Row 15
$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = "@"
$c.Value = '2.625.78'
$c.Style = "Normal"
$c = $ws.Cells.Item(15, 5)
$c.NumberFormat = "@"
$c.Value = '  +3.61%  '
$c.Style = "Normal"

# Row 16
$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = "@"
$c.Value = '15.31'
$c.Style = "Normal"
$c = $ws.Cells.Item(16, 5)
$c.NumberFormat = "@"
$c.Value = '  +0.43%  '
$c.Style = "Normal"

# Row 18
$c = $ws.Cells.Item(18, 4)
$c.NumberFormat = "@"
$c.Value = '43.304.33'
$c.Style = "Normal"
$c = $ws.Cells.Item(18, 5)
$c.NumberFormat = "@"
$c.Value = '  +1.44%  '
$c.Style = "Normal"

# Row 19
$c = $ws.Cells.Item(19, 5)
$c.NumberFormat = "@"
$c.Value = '  +3.11%  '
$c.Style = "Normal"

# Row 20
$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = "@"
$c.Value = '12.84'
$c.Style = "Normal"
$c = $ws.Cells.Item(20, 5)
$c.NumberFormat = "@"
$c.Value = '  -0.62%  '
$c.Style = "Normal"

# Row 21
$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = "@"
$c.Value = '0.0₃0971'
$c.Style = "Normal"
$c = $ws.Cells.Item(21, 5)
$c.NumberFormat = "@"
$c.Value = '  +1.46%  '
$c.Style = "Normal"

# Row 22
$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = "@"
$c.Value = '70.02'
$c.Style = "Normal"
$c = $ws.Cells.Item(22, 5)
$c.NumberFormat = "@"
$c.Value = '  +1.00%  '
$c.Style = "Normal"

# Row 23
$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = "@"
$c.Value = '255.25'
$c.Style = "Normal"
$c = $ws.Cells.Item(23, 5)
$c.NumberFormat = "@"
$c.Value = '  +1.88%  '
$c.Style = "Normal"

# Row 24
$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = "@"
$c.Value = '2.97'
$c.Style = "Normal"
$c = $ws.Cells.Item(24, 5)
$c.NumberFormat = "@"
$c.Value = '  +0.94%  '
$c.Style = "Normal"

# Row 25
$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = "@"
$c.Value = '2.11'
$c.Style = "Normal"
$c = $ws.Cells.Item(25, 5)
$c.NumberFormat = "@"
$c.Value = '  +4.73%  '
$c.Style = "Normal"

# Row 26
$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = "@"
$c.Value = '27.36'
$c.Style = "Normal"
$c = $ws.Cells.Item(26, 5)
$c.NumberFormat = "@"
$c.Value = '  +2.34%  '
$c.Style = "Normal"

# Row 27
$c = $ws.Cells.Item(27, 5)
$c.NumberFormat = "@"
$c.Value = '  -0.06%  '
$c.Style = "Normal"

# Row 28
$c = $ws.Cells.Item(28, 5)
$c.NumberFormat = "@"
$c.Value = '  -0.09%  '
$c.Style = "Normal"

# Row 29
$c = $ws.Cells.Item(29, 5)
$c.NumberFormat = "@"
$c.Value = '  +0.69%  '
$c.Style = "Normal"

# Row 30
$c = $ws.Cells.Item(30, 4)
$c.NumberFormat = "@"
$c.Value = '10.36'
$c.Style = "Normal"
$c = $ws.Cells.Item(30, 5)
$c.NumberFormat = "@"
$c.Value = '  +1.15%  '
$c.Style = "Normal"

# Row 31
$c = $ws.Cells.Item(31, 5)
$c.NumberFormat = "@"
$c.Value = '  -0.06%  '
$c.Style = "Normal"

# Row 32
$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = "@"
$c.Value = '156.56'
$c.Style = "Normal"
$c = $ws.Cells.Item(32, 5)
$c.NumberFormat = "@"
$c.Value = '  -0.85%  '
$c.Style = "Normal"

# Row 33
$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = "@"
$c.Value = '3.46'
$c.Style = "Normal"
$c = $ws.Cells.Item(33, 5)
$c.NumberFormat = "@"
$c.Value = '  +5.94%  '
$c.Style = "Normal"

# Row 34
$c = $ws.Cells.Item(34, 5)
$c.NumberFormat = "@"
$c.Value = '  +2.79%  '
$c.Style = "Normal"

# Row 35
$c = $ws.Cells.Item(35, 4)
$c.NumberFormat = "@"
$c.Value = '0.0813'
$c.Style = "Normal"
$c = $ws.Cells.Item(35, 5)
$c.NumberFormat = "@"
$c.Value = '  +4.35%  '
$c.Style = "Normal"

# Row 36
$c = $ws.Cells.Item(36, 4)
$c.NumberFormat = "@"
$c.Value = '2.71'
$c.Style = "Normal"
$c = $ws.Cells.Item(36, 5)
$c.NumberFormat = "@"
$c.Value = '  +3.21%  '
$c.Style = "Normal"

# Row 37
$c = $ws.Cells.Item(37, 4)
$c.NumberFormat = "@"
$c.Value = '18.91'
$c.Style = "Normal"
$c = $ws.Cells.Item(37, 5)
$c.NumberFormat = "@"
$c.Value = '  -2.12%  '
$c.Style = "Normal"

# Row 38
$c = $ws.Cells.Item(38, 5)
$c.NumberFormat = "@"
$c.Value = '  +1.86%  '
$c.Style = "Normal"

# Row 39
$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = "@"
$c.Value = '2.52'
$c.Style = "Normal"
$c = $ws.Cells.Item(39, 5)
$c.NumberFormat = "@"
$c.Value = '  +9.73%  '
$c.Style = "Normal"

# Row 40
$c = $ws.Cells.Item(40, 5)
$c.NumberFormat = "@"
$c.Value = '  +0.94%  '
$c.Style = "Normal"

# Row 41
$c = $ws.Cells.Item(41, 5)
$c.NumberFormat = "@"
$c.Value = '  -2.38%  '
$c.Style = "Normal"

# Row 42
$c = $ws.Cells.Item(42, 5)
$c.NumberFormat = "@"
$c.Value = '  +6.72%  '
$c.Style = "Normal"

# Row 43
$c = $ws.Cells.Item(43, 5)
$c.NumberFormat = "@"
$c.Value = '  +1.20%  '
$c.Style = "Normal"

# Row 44
$c = $ws.Cells.Item(44, 5)
$c.NumberFormat = "@"
$c.Value = '  -0.02%  '
$c.Style = "Normal"

# Row 45
$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = "@"
$c.Value = '3.26'
$c.Style = "Normal"
$c = $ws.Cells.Item(45, 5)
$c.NumberFormat = "@"
$c.Value = '  -1.13%  '
$c.Style = "Normal"

# Row 46
$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = "@"
$c.Value = '2.017.20'
$c.Style = "Normal"
$c = $ws.Cells.Item(46, 5)
$c.NumberFormat = "@"
$c.Value = '  +0.30%  '
$c.Style = "Normal"

# Row 47
$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = "@"
$c.Value = '9.02'
$c.Style = "Normal"
$c = $ws.Cells.Item(47, 5)
$c.NumberFormat = "@"
$c.Value = '  +2.15%  '
$c.Style = "Normal"

# Row 48
$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = "@"
$c.Value = '2.845.02'
$c.Style = "Normal"
$c = $ws.Cells.Item(48, 5)
$c.NumberFormat = "@"
$c.Value = '  +3.23%  '
$c.Style = "Normal"

# Row 49
$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = "@"
$c.Value = '83.78'
$c.Style = "Normal"
$c = $ws.Cells.Item(49, 5)
$c.NumberFormat = "@"
$c.Value = '  -2.05%  '
$c.Style = "Normal"

# Row 50
$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = "@"
$c.Value = '0.197'
$c.Style = "Normal"
$c = $ws.Cells.Item(50, 5)
$c.NumberFormat = "@"
$c.Value = '  +5.05%  '
$c.Style = "Normal"

# Row 51
$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = "@"
$c.Value = '74.77'
$c.Style = "Normal"
$c = $ws.Cells.Item(51, 5)
$c.NumberFormat = "@"
$c.Value = '  +1.25%  '
$c.Style = "Normal"
